$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.876.67'
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").Value = '1.629.68'
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("E6").Value = '  -1.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.28'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.31%  '
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E10").Value = '  -0.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0881'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").Value = '1.861.89'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").Value = '1.638.04'
$ws.Range("E13").Value = '  +0.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.02'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.75%  '
$ws.Range("E15").Value = '  -0.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").Value = '27.908.37'
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '228.00'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").Value = '0.0₃0719'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.92'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.35%  '
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '155.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("E29").Value = '  -0.89%  '
$ws.Range("E30").Value = '  -0.50%  '
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("D33").Value = '1.414.62'
$ws.Range("E33").Value = '  +1.28%  '
$ws.Range("E34").Value = '  +1.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.62'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.68%  '
$ws.Range("E37").Value = '  -1.31%  '
$ws.Range("E38").Value = '  -1.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.553'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("E40").Value = '  -1.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '65.86'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.94%  '
$ws.Range("E43").Value = '  -0.80%  '
$ws.Range("E44").Value = '  -0.65%  '
$ws.Range("D45").Value = '1.770.39'
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("E46").Value = '  -3.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.65'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.67%  '
$ws.Range("E48").Value = '  +1.05%  '
$ws.Range("E50").Value = '  +1.27%  '
$ws.Range("E51").Value = '  -0.13%  '
